$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"

$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 4).Value = 44706

$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100103
$ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value = 100103002
$ws.Cells.Item($row, 10).Value = "Ciruela"
$ws.Cells.Item($row, 11).Value = "Angeleno"
$ws.Cells.Item($row, 12).Value = "Segunda"
$ws.Cells.Item($row, 13).Value = 300
$ws.Cells.Item($row, 14).Value = 15000
$ws.Cells.Item($row, 15).Value = 16000
$ws.Cells.Item($row, 16).Value = 15500
$ws.Cells.Item($row, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 861
$ws.Cells.Item($row, 20).Value = 18
